$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old first data row (row 2): this shifts all subsequent rows
# up by one, matching the new date/year alignment (A, B, C, D columns).
$ws.Rows(2).Delete()

# Clear any stale leftover in column E from the old row 19 that would
# otherwise remain after the shift (defensive; also ensures E1:E5 blank).
$ws.Range("E2:E5").ClearContents()

# Recalculated y_1_forecast (column E) values for rows 6-18.
$newE = @{
    6  = 0.9288717675470126
    7  = 1.416624765035412
    8  = 1.075154359849861
    9  = 1.214249019249602
    10 = 1.180565832117297
    11 = 1.638669199130427
    12 = 1.540918326052476
    13 = 0.5773070399857971
    14 = -0.415982961498651
    15 = 1.426719405738508
    16 = 0.5205511175203181
    17 = 0.3517304536567734
    18 = 0.4186921370205043
}

foreach ($row in $newE.Keys) {
    $ws.Cells.Item($row, 5).Value = $newE[$row]
}
